{"js": "// Replace the 100 arithmetic-problem cells in the single 20x5 table with\n// their updated values, looked up by the cell's current (old) text so the\n// mapping is robust to any re-ordering and independent of row/col indices.\nconst replacements = {\"49-11=\": \"46-8=\", \"21-1=\": \"47+45=\", \"92-7=\": \"72-36=\", \"51+4=\": \"21+65=\", \"78-27=\": \"44+8=\", \"75+13=\": \"25+36=\", \"0+44=\": \"25+61=\", \"26+13=\": \"45-15=\", \"46+3=\": \"52-10=\", \"93-55=\": \"28+48=\", \"83-27=\": \"69-6=\", \"68-48=\": \"96-90=\", \"85-17=\": \"30+48=\", \"34+50=\": \"25+70=\", \"88-80=\": \"2+59=\", \"21+9=\": \"63-3=\", \"83-8=\": \"91-29=\", \"46-35=\": \"39-28=\", \"0+52=\": \"7+6=\", \"46+24=\": \"14-0=\", \"66-64=\": \"86-3=\", \"64-7=\": \"53+18=\", \"3+32=\": \"12+57=\", \"89-63=\": \"98-56=\", \"62+12=\": \"21+26=\", \"4+59=\": \"65+24=\", \"62+8=\": \"19+34=\", \"78+13=\": \"58+37=\", \"48-17=\": \"81-60=\", \"36-19=\": \"63-7=\", \"86-46=\": \"97-46=\", \"0+57=\": \"42+16=\", \"14+14=\": \"84-10=\", \"14-7=\": \"91-62=\", \"8+45=\": \"63-10=\", \"67-23=\": \"54-15=\", \"77-73=\": \"57-20=\", \"77-24=\": \"86-29=\", \"72-65=\": \"59+39=\", \"18+14=\": \"19+23=\", \"16+59=\": \"47+20=\", \"49-42=\": \"9+26=\", \"81-72=\": \"22+53=\", \"72+23=\": \"87-17=\", \"37+23=\": \"67-41=\", \"78-70=\": \"8+10=\", \"16+41=\": \"57-22=\", \"78-52=\": \"83+5=\", \"36-9=\": \"11+51=\", \"53+13=\": \"16+53=\", \"35+56=\": \"52-4=\", \"8+67=\": \"55+42=\", \"52-28=\": \"59-4=\", \"18+22=\": \"2+35=\", \"40+51=\": \"79-36=\", \"5+49=\": \"84+4=\", \"24+72=\": \"30-19=\", \"6+84=\": \"34+32=\", \"84-25=\": \"2+5=\", \"47-39=\": \"27-26=\", \"55-54=\": \"3+21=\", \"9+78=\": \"1+27=\", \"5+91=\": \"16+62=\", \"87-21=\": \"82-62=\", \"87+2=\": \"77-10=\", \"94-67=\": \"0+88=\", \"14+6=\": \"93-57=\", \"15+39=\": \"66-25=\", \"2+49=\": \"57+32=\", \"16+60=\": \"48-16=\", \"75-31=\": \"86-30=\", \"31+67=\": \"25+19=\", \"12+70=\": \"43+40=\", \"73-39=\": \"66-38=\", \"16-9=\": \"57-13=\", \"37+15=\": \"23+54=\", \"9+43=\": \"20+76=\", \"15-5=\": \"94-16=\", \"45+15=\": \"97-62=\", \"0+16=\": \"67+1=\", \"46-37=\": \"12+14=\", \"21+55=\": \"15+7=\", \"6+25=\": \"69+7=\", \"48-2=\": \"82-71=\", \"87-32=\": \"78+6=\", \"75-21=\": \"63-46=\", \"97-49=\": \"7+91=\", \"70+13=\": \"89-29=\", \"84-49=\": \"40+2=\", \"33+48=\": \"83-69=\", \"65-38=\": \"75+7=\", \"13+19=\": \"35+9=\", \"9+22=\": \"80-71=\", \"76-0=\": \"47-25=\", \"24+21=\": \"37-16=\", \"66-18=\": \"42+21=\", \"80-16=\": \"10+71=\", \"20+30=\": \"15+33=\", \"95+4=\": \"82-57=\", \"25+1=\": \"0+80=\"};\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values,rowCount,columnCount\");\nawait context.sync();\n\nconst oldValues = table.values;\nconst newValues = oldValues.map(row =>\n  row.map(cellText => {\n    const trimmed = cellText.trim();\n    return Object.prototype.hasOwnProperty.call(replacements, trimmed)\n      ? replacements[trimmed]\n      : cellText;\n  })\n);\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem cells in the single 20x5 table with\n# their updated values, looked up by the cell's current (old) text so the\n# mapping is robust regardless of row/col indices.\n\n$replacements = @{\n  '49-11=' = '46-8='\n  '21-1=' = '47+45='\n  '92-7=' = '72-36='\n  '51+4=' = '21+65='\n  '78-27=' = '44+8='\n  '75+13=' = '25+36='\n  '0+44=' = '25+61='\n  '26+13=' = '45-15='\n  '46+3=' = '52-10='\n  '93-55=' = '28+48='\n  '83-27=' = '69-6='\n  '68-48=' = '96-90='\n  '85-17=' = '30+48='\n  '34+50=' = '25+70='\n  '88-80=' = '2+59='\n  '21+9=' = '63-3='\n  '83-8=' = '91-29='\n  '46-35=' = '39-28='\n  '0+52=' = '7+6='\n  '46+24=' = '14-0='\n  '66-64=' = '86-3='\n  '64-7=' = '53+18='\n  '3+32=' = '12+57='\n  '89-63=' = '98-56='\n  '62+12=' = '21+26='\n  '4+59=' = '65+24='\n  '62+8=' = '19+34='\n  '78+13=' = '58+37='\n  '48-17=' = '81-60='\n  '36-19=' = '63-7='\n  '86-46=' = '97-46='\n  '0+57=' = '42+16='\n  '14+14=' = '84-10='\n  '14-7=' = '91-62='\n  '8+45=' = '63-10='\n  '67-23=' = '54-15='\n  '77-73=' = '57-20='\n  '77-24=' = '86-29='\n  '72-65=' = '59+39='\n  '18+14=' = '19+23='\n  '16+59=' = '47+20='\n  '49-42=' = '9+26='\n  '81-72=' = '22+53='\n  '72+23=' = '87-17='\n  '37+23=' = '67-41='\n  '78-70=' = '8+10='\n  '16+41=' = '57-22='\n  '78-52=' = '83+5='\n  '36-9=' = '11+51='\n  '53+13=' = '16+53='\n  '35+56=' = '52-4='\n  '8+67=' = '55+42='\n  '52-28=' = '59-4='\n  '18+22=' = '2+35='\n  '40+51=' = '79-36='\n  '5+49=' = '84+4='\n  '24+72=' = '30-19='\n  '6+84=' = '34+32='\n  '84-25=' = '2+5='\n  '47-39=' = '27-26='\n  '55-54=' = '3+21='\n  '9+78=' = '1+27='\n  '5+91=' = '16+62='\n  '87-21=' = '82-62='\n  '87+2=' = '77-10='\n  '94-67=' = '0+88='\n  '14+6=' = '93-57='\n  '15+39=' = '66-25='\n  '2+49=' = '57+32='\n  '16+60=' = '48-16='\n  '75-31=' = '86-30='\n  '31+67=' = '25+19='\n  '12+70=' = '43+40='\n  '73-39=' = '66-38='\n  '16-9=' = '57-13='\n  '37+15=' = '23+54='\n  '9+43=' = '20+76='\n  '15-5=' = '94-16='\n  '45+15=' = '97-62='\n  '0+16=' = '67+1='\n  '46-37=' = '12+14='\n  '21+55=' = '15+7='\n  '6+25=' = '69+7='\n  '48-2=' = '82-71='\n  '87-32=' = '78+6='\n  '75-21=' = '63-46='\n  '97-49=' = '7+91='\n  '70+13=' = '89-29='\n  '84-49=' = '40+2='\n  '33+48=' = '83-69='\n  '65-38=' = '75+7='\n  '13+19=' = '35+9='\n  '9+22=' = '80-71='\n  '76-0=' = '47-25='\n  '24+21=' = '37-16='\n  '66-18=' = '42+21='\n  '80-16=' = '10+71='\n  '20+30=' = '15+33='\n  '95+4=' = '82-57='\n  '25+1=' = '0+80='\n}\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $range = $cell.Range\n    $txt = $range.Text\n    # Cell text ends with a cell-mark (chr 7) preceded by a paragraph mark (chr 13);\n    # strip those control characters before comparing/looking up.\n    $clean = $txt.TrimEnd([char]13, [char]7)\n    if ($replacements.ContainsKey($clean)) {\n      $cell.Range.Text = $replacements[$clean]\n    }\n  }\n}\n"}
